# Fix AV rappel issue
# Rebuild the data rows (2-9) of the "Etat Taxes" sheet with the corrected
# contract / tenant information and rappel amounts, adding the new rows
# that were missing and moving the totals row down to row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the 15 column values (A..O) for one row.
$rows = @(
    ,@("555/RRR/AV3", "Direction régionale", "B171710",  "NADIA BADRANE",    "non", "mensuelle", 0, "--",  0, "--", 0, 16000, 0,    "--", 16000)
    ,@("555/RRR/AV3", "Direction régionale", "IB43905",  "NHILA BELGACEM",   "non", "mensuelle", 0, "--",  0, "--", 0, 16000, 0,    "--", 16000)
    ,@("555/RRR/AV3", "Direction régionale", "B171710",  "NADIA BADRANE",    "non", "mensuelle", 0, "--",  0, "--", 0, 2000,  0,    "--", 2000)
    ,@("555/RRR/AV3", "Direction régionale", "IB43905",  "NHILA BELGACEM",   "non", "mensuelle", 0, "--",  0, "--", 0, 2000,  0,    "--", 2000)
    ,@("555/RRR/AV3", "Direction régionale", "B171710",  "NADIA BADRANE",    "non", "mensuelle", 0, 1000,  0, 0,    0, 0,     0,    "--", 1000)
    ,@("555/RRR/AV3", "Direction régionale", "IB43905",  "NHILA BELGACEM",   "non", "mensuelle", 0, 1000,  0, 0,    0, 0,     0,    "--", 1000)
    ,@("000/CCCC",    "Direction régionale", "BK646476", "DOUNIA LAMKADDAM", "non", "mensuelle", 0, 2000,  0, 0,    0, 0,     0,    "--", 2000)
    ,@(" ", " ", " ", " ", " ", " ", " ", 4000, 0, 0, 0, 36000, 0, 0, 40000)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

Write-Host "Updated rows 2-9 with corrected AV rappel data"
